# Conform date text in column A to a standard "MM.DD.YY" form:
# "12.5.18" -> "12.05.18"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq "12.5.18") {
        # Force the new value to stay literal text (Excel would otherwise
        # reinterpret a date-looking string like "12.05.18" as a date
        # serial number), then restore the default "Normal" style so the
        # cell keeps no explicit style override.
        $cell.NumberFormat = "@"
        $cell.Value2 = "12.05.18"
        $cell.Style = "Normal"
    }
}

# Move the saved selection/active cell from A6:XFD6 to A4
$ws.Range("A4").Select()
